$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value2 = 52882.09
$ws.Range("I51").Value2 = 8152.6665
$ws.Range("K51").Value2 = 8152.6665
$ws.Range("M51").Value2 = -7668.6665

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value2 = 76123.25
$ws.Range("I100").Value2 = 84417.086
$ws.Range("K100").Value2 = 84417.086
$ws.Range("M100").Value2 = -83876.086

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value2 = 9404361
$ws.Range("I116").Value2 = 12224520
$ws.Range("K116").Value2 = 12224520
$ws.Range("M116").Value2 = -12221078

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value2 = 1541444.5
$ws.Range("I132").Value2 = 3101.1155
$ws.Range("K132").Value2 = 9303.3465
$ws.Range("M132").Value2 = -6773.3465

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value2 = 196928.44
$ws.Range("I138").Value2 = 823882.1
$ws.Range("J138").Value2 = 4019.6155
$ws.Range("K138").Value2 = 2471646.3
$ws.Range("L138").Value2 = 12058.8465
$ws.Range("M138").Value2 = -2466506.3
$ws.Range("N138").Value2 = -22338.8465

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value2 = 5197.7334
$ws.Range("I2").Value2 = 9503.333000000001
$ws.Range("J2").Value2 = 2327.3333
$ws.Range("K2").Value2 = 9503.333000000001
$ws.Range("L2").Value2 = 2327.3333
$ws.Range("M2").Value2 = -9390.333000000001
$ws.Range("N2").Value2 = -2553.3333

# Row 41: Skillet Scandal / White Skillet
$ws.Range("H41").Value2 = 3616
$ws.Range("I41").Value2 = 2645
$ws.Range("K41").Value2 = 2645
$ws.Range("M41").Value2 = -2231

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value2 = 9087
$ws.Range("I45").Value2 = 8130.5
$ws.Range("K45").Value2 = 8130.5
$ws.Range("M45").Value2 = -7753.5

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value2 = 1676.9259
$ws.Range("I74").Value2 = 635.5454999999999
$ws.Range("K74").Value2 = 635.5454999999999
$ws.Range("M74").Value2 = 238.4545000000001

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value2 = 1676.9259
$ws.Range("I77").Value2 = 635.5454999999999
$ws.Range("K77").Value2 = 3177.7275
$ws.Range("M77").Value2 = 1190.2725

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value2 = 5197.7334
$ws.Range("I116").Value2 = 9503.333000000001
$ws.Range("J116").Value2 = 2327.3333
$ws.Range("K116").Value2 = 9503.333000000001
$ws.Range("L116").Value2 = 2327.3333
$ws.Range("M116").Value2 = -7209.333000000001
$ws.Range("N116").Value2 = -6915.3333

# Row 121: Shield to Shield / Dwarven Mythril Shield
$ws.Range("H121").Value2 = 0
$ws.Range("J121").Value2 = 0
$ws.Range("L121").Value2 = 0
$ws.Range("N121").ClearContents()

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value2 = 2612.8125
$ws.Range("I132").Value2 = 766.5714
$ws.Range("K132").Value2 = 2299.7142
$ws.Range("M132").Value2 = 230.2857999999997

# Row 135: Forgiveness for My Shins / Ruthenium Sabatons of Fending
$ws.Range("H135").Value2 = 83441.78
$ws.Range("J135").Value2 = 83441.78
$ws.Range("L135").Value2 = 83441.78
$ws.Range("N135").Value2 = -93581.78

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value2 = 5197.7334
$ws.Range("I3").Value2 = 9503.333000000001
$ws.Range("J3").Value2 = 2327.3333
$ws.Range("K3").Value2 = 9503.333000000001
$ws.Range("L3").Value2 = 2327.3333
$ws.Range("M3").Value2 = -9389.333000000001
$ws.Range("N3").Value2 = -2555.3333

# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value2 = 247.8
$ws.Range("I22").Value2 = 224.75
$ws.Range("J22").Value2 = 340
$ws.Range("K22").Value2 = 224.75
$ws.Range("L22").Value2 = 340
$ws.Range("M22").Value2 = -51.75
$ws.Range("N22").Value2 = -686

# Row 52: File That under Whatever / Mythril File
$ws.Range("H52").Value2 = 52438.777
$ws.Range("J52").Value2 = 55487.25
$ws.Range("L52").Value2 = 55487.25
$ws.Range("N52").Value2 = -56013.25

# Row 121: Keeping Loyalty / Dwarven Mythril File
$ws.Range("H121").Value2 = 52438.777
$ws.Range("J121").Value2 = 55487.25
$ws.Range("L121").Value2 = 55487.25
$ws.Range("N121").Value2 = -58981.25

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value2 = 6110.7188
$ws.Range("I134").Value2 = 6637.5
$ws.Range("K134").Value2 = 19912.5
$ws.Range("M134").Value2 = -17377.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value2 = 8972.576999999999
$ws.Range("I31").Value2 = 11035.706
$ws.Range("J31").Value2 = 5075.5557
$ws.Range("K31").Value2 = 11035.706
$ws.Range("L31").Value2 = 5075.5557
$ws.Range("M31").Value2 = -10740.706
$ws.Range("N31").Value2 = -5665.5557

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value2 = 8972.576999999999
$ws.Range("I34").Value2 = 11035.706
$ws.Range("J34").Value2 = 5075.5557
$ws.Range("K34").Value2 = 11035.706
$ws.Range("L34").Value2 = 5075.5557
$ws.Range("M34").Value2 = -10833.706
$ws.Range("N34").Value2 = -5479.5557

# Row 108: Just Starting Out / White Oak Fishing Rod
$ws.Range("H108").Value2 = 44097.668
$ws.Range("J108").Value2 = 44097.668
$ws.Range("L108").Value2 = 44097.668
$ws.Range("N108").Value2 = -51777.668

# Row 109: Playing the Market / White Oak Necklace
$ws.Range("H109").Value2 = 52000
$ws.Range("J109").Value2 = 52000
$ws.Range("L109").Value2 = 52000
$ws.Range("N109").Value2 = -54080

# Row 114: Ground to a Halt / White Ash Grinding Wheel
$ws.Range("H114").Value2 = 21277.5
$ws.Range("J114").Value2 = 21277.5
$ws.Range("L114").Value2 = 21277.5
$ws.Range("N114").Value2 = -29955.5

# Row 117: Sleep on It / Sandteak Spinning Wheel
$ws.Range("H117").Value2 = 65555
$ws.Range("J117").Value2 = 65555
$ws.Range("L117").Value2 = 65555
$ws.Range("N117").Value2 = -74733

# Row 131: An Integral Reward / Integral Necklace of Crafting
$ws.Range("H131").Value2 = 0
$ws.Range("J131").Value2 = 0
$ws.Range("L131").Value2 = 0
$ws.Range("N131").ClearContents()

# Row 141: No Greater Treasure / Claro Walnut Necklace of Gathering
$ws.Range("H141").Value2 = 171413.97
$ws.Range("J141").Value2 = 183870.25
$ws.Range("L141").Value2 = 183870.25
$ws.Range("N141").Value2 = -194230.25

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water / Boiled Egg
$ws.Range("H4").Value2 = 78891010
$ws.Range("I4").Value2 = 71206870
$ws.Range("K4").Value2 = 213620610
$ws.Range("M4").Value2 = -213620498

# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value2 = 239458.92
$ws.Range("I5").Value2 = 454.21054
$ws.Range("J5").Value2 = 436897.6
$ws.Range("K5").Value2 = 1362.63162
$ws.Range("L5").Value2 = 1310692.8
$ws.Range("M5").Value2 = -1250.63162
$ws.Range("N5").Value2 = -1310916.8

# Row 33: Cooking with Gas / Chicken Stock
$ws.Range("H33").Value2 = 408
$ws.Range("I33").Value2 = 441.33334
$ws.Range("J33").Value2 = 391.33334
$ws.Range("K33").Value2 = 2648.00004
$ws.Range("L33").Value2 = 2348.00004
$ws.Range("M33").Value2 = -2365.00004
$ws.Range("N33").Value2 = -2914.00004

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value2 = 55561100
$ws.Range("J131").Value2 = 2190.1428
$ws.Range("L131").Value2 = 6570.428400000001
$ws.Range("N131").Value2 = -16650.4284

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value2 = 239458.92
$ws.Range("I135").Value2 = 454.21054
$ws.Range("J135").Value2 = 436897.6
$ws.Range("K135").Value2 = 4087.89486
$ws.Range("L135").Value2 = 3932078.4
$ws.Range("M135").Value2 = -1552.89486
$ws.Range("N135").Value2 = -3937148.4

$ws = $wb.Worksheets.Item("GSM")
# Row 19: Better Four Eyes than None / Brass Spectacles
$ws.Range("H19").Value2 = 1527
$ws.Range("I19").Value2 = 1527
$ws.Range("K19").Value2 = 1527
$ws.Range("M19").Value2 = -1239

# Row 62: The Goggles, They Do Naught / Mythrite Goggles of Gathering
$ws.Range("H62").Value2 = 29000
$ws.Range("I62").Value2 = 29000
$ws.Range("K62").Value2 = 29000
$ws.Range("M62").Value2 = -28314

# Row 65: Peril Never Wore Safety Goggles (L) / Mythrite Goggles of Gathering
$ws.Range("H65").Value2 = 29000
$ws.Range("I65").Value2 = 29000
$ws.Range("K65").Value2 = 87000
$ws.Range("M65").Value2 = -83568

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value2 = 6731.913
$ws.Range("I113").Value2 = 9062.333000000001
$ws.Range("K113").Value2 = 9062.333000000001
$ws.Range("M113").Value2 = -6892.333000000001

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value2 = 2042.9048
$ws.Range("I132").Value2 = 1997.9395
$ws.Range("K132").Value2 = 5993.818499999999
$ws.Range("M132").Value2 = -3463.818499999999

$ws = $wb.Worksheets.Item("LTW")
# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value2 = 1428.7333
$ws.Range("I55").Value2 = 509
$ws.Range("K55").Value2 = 509
$ws.Range("M55").Value2 = -336

# Row 63: From Mud to Mourning / Archaeoskin Jackboots of Gathering
$ws.Range("H63").Value2 = 44998.5
$ws.Range("J63").Value2 = 44998.5
$ws.Range("L63").Value2 = 44998.5
$ws.Range("N63").Value2 = -46496.5

# Row 66: These Boots Are Made for Hawkin' (L) / Archaeoskin Jackboots of Gathering
$ws.Range("H66").Value2 = 44998.5
$ws.Range("J66").Value2 = 44998.5
$ws.Range("L66").Value2 = 134995.5
$ws.Range("N66").Value2 = -142483.5

# Row 109: Band Substances / Smilodonskin Wristband
$ws.Range("H109").Value2 = 73983.5
$ws.Range("J109").Value2 = 73983.5
$ws.Range("L109").Value2 = 73983.5
$ws.Range("N109").Value2 = -76757.5

# Row 134: Freezing Fingers / Crocodileskin Fingerless Gloves of Striking
$ws.Range("H134").Value2 = 66249.75
$ws.Range("J134").Value2 = 66249.75
$ws.Range("L134").Value2 = 66249.75
$ws.Range("N134").Value2 = -76389.75

$ws = $wb.Worksheets.Item("WVR")
# Row 21: Don't Trew So Hard / Initiate's Slops
$ws.Range("H21").Value2 = 39000
$ws.Range("J21").Value2 = 0
$ws.Range("L21").Value2 = 0
$ws.Range("N21").ClearContents()

# Row 35: Pantser Corps / Initiate's Slops
$ws.Range("H35").Value2 = 39000
$ws.Range("J35").Value2 = 0
$ws.Range("L35").Value2 = 0
$ws.Range("N35").ClearContents()

# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value2 = 46128.57
$ws.Range("I107").Value2 = 8500
$ws.Range("J107").Value2 = 61180
$ws.Range("K107").Value2 = 25500
$ws.Range("L107").Value2 = 183540
$ws.Range("M107").Value2 = -23580
$ws.Range("N107").Value2 = -187380

# Row 133: Begin with the Basics / Snow Cotton Jacket
$ws.Range("H133").Value2 = 70000
$ws.Range("J133").Value2 = 70000
$ws.Range("L133").Value2 = 70000
$ws.Range("N133").Value2 = -80120

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value2 = 533439.2
$ws.Range("I136").Value2 = 617913.3
$ws.Range("K136").Value2 = 1853739.9
$ws.Range("M136").Value2 = -1851189.9

